# Weekly Fruta/Hortaliza update: a new week of price data was published.
# Three new daily observations were inserted at the top of the recent data
# block (rows 192-194), pushing all the previously-existing rows 192-217
# down to rows 195-220.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the three new rows; everything currently at rows 192:217
# (and their formatting) slides down to 195:220.
$ws.Rows("192:194").Insert()

# Row 192 - new observation, "Primera" quality, Provincia de Chacabuco
$ws.Cells.Item(192, 1).Value = 9
$ws.Cells.Item(192, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(192, 3).Value = "Metropolitana"
$ws.Cells.Item(192, 4).Value = 44504
$ws.Cells.Item(192, 5).Value = 13
$ws.Cells.Item(192, 6).Value = 100112052
$ws.Cells.Item(192, 7).Value = "Albahaca"
$ws.Cells.Item(192, 8).Value = "Sin especificar"
$ws.Cells.Item(192, 9).Value = "Primera"
$ws.Cells.Item(192, 10).Value = 52
$ws.Cells.Item(192, 11).Value = 6000
$ws.Cells.Item(192, 12).Value = 7000
$ws.Cells.Item(192, 13).Value = 6500
$ws.Cells.Item(192, 14).Value = "`$/docena de matas"
$ws.Cells.Item(192, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(192, 16).Value = 1083
$ws.Cells.Item(192, 17).Value = 6
$ws.Cells.Item(192, 18).Value = "Hortaliza"

# Row 193 - new observation, "Primera" quality, Region de Arica y Parinacota
$ws.Cells.Item(193, 1).Value = 9
$ws.Cells.Item(193, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(193, 3).Value = "Metropolitana"
$ws.Cells.Item(193, 4).Value = 44504
$ws.Cells.Item(193, 5).Value = 13
$ws.Cells.Item(193, 6).Value = 100112052
$ws.Cells.Item(193, 7).Value = "Albahaca"
$ws.Cells.Item(193, 8).Value = "Sin especificar"
$ws.Cells.Item(193, 9).Value = "Primera"
$ws.Cells.Item(193, 10).Value = 160
$ws.Cells.Item(193, 11).Value = 5000
$ws.Cells.Item(193, 12).Value = 5000
$ws.Cells.Item(193, 13).Value = 5000
$ws.Cells.Item(193, 14).Value = "`$/paquete"
$ws.Cells.Item(193, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(193, 16).Value = 5000
$ws.Cells.Item(193, 17).Value = 1
$ws.Cells.Item(193, 18).Value = "Hortaliza"

# Row 194 - new observation, "Segunda" quality, Region de Arica y Parinacota
$ws.Cells.Item(194, 1).Value = 9
$ws.Cells.Item(194, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(194, 3).Value = "Metropolitana"
$ws.Cells.Item(194, 4).Value = 44504
$ws.Cells.Item(194, 5).Value = 13
$ws.Cells.Item(194, 6).Value = 100112052
$ws.Cells.Item(194, 7).Value = "Albahaca"
$ws.Cells.Item(194, 8).Value = "Sin especificar"
$ws.Cells.Item(194, 9).Value = "Segunda"
$ws.Cells.Item(194, 10).Value = 61
$ws.Cells.Item(194, 11).Value = 4000
$ws.Cells.Item(194, 12).Value = 4000
$ws.Cells.Item(194, 13).Value = 4000
$ws.Cells.Item(194, 14).Value = "`$/paquete"
$ws.Cells.Item(194, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(194, 16).Value = 4000
$ws.Cells.Item(194, 17).Value = 1
$ws.Cells.Item(194, 18).Value = "Hortaliza"
